$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = "Bidar"
$ws.Range("G28").Value = "Bidar"
$ws.Range("G31").Value = "Bidar"
$ws.Range("G52").Value = "Ballari (Bellary)"
$ws.Range("G60").Value = "Bidar"
$ws.Range("G73").Value = "Bidar"
